$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $style
}

Set-TextValue $ws.Range("D2") '307.82'
Set-TextValue $ws.Range("E2") '-4.33%'
Set-TextValue $ws.Range("G2") '15'

Set-TextValue $ws.Range("D3") '39.83'
Set-TextValue $ws.Range("E3") '-6.69%'
Set-TextValue $ws.Range("G3") '15'

Set-TextValue $ws.Range("D4") '5.086'
Set-TextValue $ws.Range("E4") '-2.46%'
Set-TextValue $ws.Range("G4") '15'

Set-TextValue $ws.Range("D5") '0.07694'
Set-TextValue $ws.Range("E5") '-6.38%'
Set-TextValue $ws.Range("G5") '15'

Set-TextValue $ws.Range("D6") '4.244'
Set-TextValue $ws.Range("E6") '-1.96%'
Set-TextValue $ws.Range("G6") '15'

Set-TextValue $ws.Range("D7") '1.626'
Set-TextValue $ws.Range("E7") '-8.31%'
Set-TextValue $ws.Range("G7") '15'

Set-TextValue $ws.Range("D8") '0.9189'
Set-TextValue $ws.Range("E8") '-3.32%'
Set-TextValue $ws.Range("G8") '15'

Set-TextValue $ws.Range("D9") '0.1036'
Set-TextValue $ws.Range("E9") '-7.62%'
Set-TextValue $ws.Range("G9") '15'

Set-TextValue $ws.Range("D10") '0.1786'
Set-TextValue $ws.Range("E10") '-5.50%'
Set-TextValue $ws.Range("G10") '15'

Set-TextValue $ws.Range("D11") '0.09284'
Set-TextValue $ws.Range("E11") '-1.45%'
Set-TextValue $ws.Range("G11") '15'

Set-TextValue $ws.Range("D12") '0.04436'
Set-TextValue $ws.Range("E12") '-4.38%'
Set-TextValue $ws.Range("G12") '15'

Set-TextValue $ws.Range("D13") '0.1054'
Set-TextValue $ws.Range("E13") '-0.43%'
Set-TextValue $ws.Range("G13") '15'

Set-TextValue $ws.Range("D14") '0.001235'
Set-TextValue $ws.Range("E14") '-4.47%'
Set-TextValue $ws.Range("G14") '15'

Set-TextValue $ws.Range("D15") '0.005868'
Set-TextValue $ws.Range("E15") '2.84%'
Set-TextValue $ws.Range("G15") '15'

Set-TextValue $ws.Range("E16") '2,410.52%'
Set-TextValue $ws.Range("G16") '15'

Set-TextValue $ws.Range("D17") '3.364'
Set-TextValue $ws.Range("E17") '0.31%'
Set-TextValue $ws.Range("G17") '15'

Set-TextValue $ws.Range("D18") '2.421'
Set-TextValue $ws.Range("E18") '-5.51%'
Set-TextValue $ws.Range("G18") '15'

Set-TextValue $ws.Range("D19") '0.3312'
Set-TextValue $ws.Range("E19") '-1.66%'
Set-TextValue $ws.Range("G19") '15'

Set-TextValue $ws.Range("D20") '6.901'
Set-TextValue $ws.Range("E20") '-7.58%'
Set-TextValue $ws.Range("G20") '15'

Set-TextValue $ws.Range("D21") '0.1346'
Set-TextValue $ws.Range("E21") '-3.17%'
Set-TextValue $ws.Range("G21") '15'

Set-TextValue $ws.Range("D22") '0.2705'
Set-TextValue $ws.Range("E22") '5.94%'
Set-TextValue $ws.Range("G22") '15'

Set-TextValue $ws.Range("D23") '0.04140'
Set-TextValue $ws.Range("E23") '-0.67%'
Set-TextValue $ws.Range("G23") '15'

Set-TextValue $ws.Range("D24") '0.001203'
Set-TextValue $ws.Range("E24") '-3.93%'
Set-TextValue $ws.Range("G24") '15'

Set-TextValue $ws.Range("D25") '0.004098'
Set-TextValue $ws.Range("E25") '-4.22%'
Set-TextValue $ws.Range("G25") '15'

Set-TextValue $ws.Range("D26") '0.0001299'
Set-TextValue $ws.Range("E26") '6.29%'
Set-TextValue $ws.Range("G26") '15'

Set-TextValue $ws.Range("G27") '15'

Set-TextValue $ws.Range("G28") '15'

Set-TextValue $ws.Range("G29") '15'

Set-TextValue $ws.Range("G30") '15'

Set-TextValue $ws.Range("G31") '15'

Set-TextValue $ws.Range("G32") '15'

Set-TextValue $ws.Range("G33") '15'

Set-TextValue $ws.Range("G34") '15'

Set-TextValue $ws.Range("G35") '15'

Set-TextValue $ws.Range("G36") '15'

Set-TextValue $ws.Range("G37") '15'

Set-TextValue $ws.Range("D38") '0.02468'
Set-TextValue $ws.Range("E38") '-7.00%'
Set-TextValue $ws.Range("G38") '15'

Set-TextValue $ws.Range("D39") '0.05189'
Set-TextValue $ws.Range("E39") '-8.21%'
Set-TextValue $ws.Range("G39") '15'

Set-TextValue $ws.Range("D40") '0.007939'
Set-TextValue $ws.Range("E40") '-2.64%'
Set-TextValue $ws.Range("G40") '15'

Set-TextValue $ws.Range("D41") '0.1320'
Set-TextValue $ws.Range("E41") '-5.86%'
Set-TextValue $ws.Range("G41") '15'

Set-TextValue $ws.Range("D42") '0.007046'
Set-TextValue $ws.Range("E42") '8.84%'
Set-TextValue $ws.Range("G42") '15'

Set-TextValue $ws.Range("D43") '0.001948'
Set-TextValue $ws.Range("E43") '-4.69%'
Set-TextValue $ws.Range("G43") '15'

Set-TextValue $ws.Range("D44") '0.007974'
Set-TextValue $ws.Range("E44") '3.48%'
Set-TextValue $ws.Range("G44") '15'

Set-TextValue $ws.Range("D45") '0.3072'
Set-TextValue $ws.Range("E45") '-12.18%'
Set-TextValue $ws.Range("G45") '15'

Set-TextValue $ws.Range("D46") '0.00006386'
Set-TextValue $ws.Range("E46") '-5.82%'
Set-TextValue $ws.Range("G46") '15'

Set-TextValue $ws.Range("D47") '0.00000000749'
Set-TextValue $ws.Range("E47") '-0.25%'
Set-TextValue $ws.Range("G47") '15'

Set-TextValue $ws.Range("D48") '0.002997'
Set-TextValue $ws.Range("E48") '-27.02%'
Set-TextValue $ws.Range("G48") '15'

Set-TextValue $ws.Range("D49") '0.004497'
Set-TextValue $ws.Range("E49") '34.41%'
Set-TextValue $ws.Range("G49") '15'

Set-TextValue $ws.Range("D50") '0.00002098'
Set-TextValue $ws.Range("E50") '-0.25%'
Set-TextValue $ws.Range("G50") '15'

Set-TextValue $ws.Range("D51") '0.0001998'
Set-TextValue $ws.Range("E51") '-0.25%'
Set-TextValue $ws.Range("G51") '15'
